$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete data rows (5 through 9); keep only header + 3 data rows
$ws.Rows("5:9").Delete()

# Fix typo in the header captions (prestation -> préstation)
$ws.Range("D1").Value = "CMP préstation unitaire"
$ws.Range("E1").Value = "Fournitures préstation unitaire"

# Row 2
$ws.Range("A2").Value = 3.1
$ws.Range("B2").Value = "Installation d'un lecteur ""passe sans contact"""
$ws.Range("C2").Value = "JOUR"
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 124.16
$ws.Range("F2").Value = 275.1
$ws.Range("G2").Value = 972.4799999999999
$ws.Range("H2").Value = 1247.58

# Row 3
$ws.Range("A3").Value = 3.1
$ws.Range("B3").Value = "Installation d'un lecteur ""passe sans contact"""
$ws.Range("C3").Value = "JOUR"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 124.16
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

# Row 4
$ws.Range("A4").Value = 3.1
$ws.Range("B4").Value = "Installation d'un lecteur ""passe sans contact"""
$ws.Range("C4").Value = "JOUR"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 124.16
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
